# Commit: Fruta / hortaliza, semanal
# A new weekly observation is inserted at row 405 (Perejil, Mercado Mayorista
# Lo Valledor de Santiago), pushing the existing rows 405:427 down to 406:428.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 405; everything from old row 405 on
# down shifts to row 406 onward (same as old row405 values now live in 406).
$ws.Rows.Item(405).Insert()

# Populate the freshly inserted row 405 with the new weekly record. The
# non-numeric / categorical columns (A,B,C,E,F,G,H,I,N,O,Q,R) keep the same
# values that were already in this slot (identical to what is now row 406),
# only the date + price/volume figures (D,J,K,L,M,P) are new.
$ws.Cells.Item(405, 1).Value = 6
$ws.Cells.Item(405, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(405, 3).Value = "Metropolitana"
$ws.Cells.Item(405, 4).Value = 44578
$ws.Cells.Item(405, 5).Value = 13
$ws.Cells.Item(405, 6).Value = 100112044
$ws.Cells.Item(405, 7).Value = "Perejil"
$ws.Cells.Item(405, 8).Value = "Sin especificar"
$ws.Cells.Item(405, 9).Value = "Primera"
$ws.Cells.Item(405, 10).Value = 150
$ws.Cells.Item(405, 11).Value = 14000
$ws.Cells.Item(405, 12).Value = 15000
$ws.Cells.Item(405, 13).Value = 14400
$ws.Cells.Item(405, 14).Value = "$/docena de atados"
$ws.Cells.Item(405, 15).Value = "Región Metropolitana"
$ws.Cells.Item(405, 16).Value = 4800
$ws.Cells.Item(405, 17).Value = 3
$ws.Cells.Item(405, 18).Value = "Hortaliza"

# Match the date-format style used by the rest of the "Fecha" column (D).
$ws.Cells.Item(405, 4).NumberFormat = $ws.Cells.Item(406, 4).NumberFormat
